$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data
$ws.Range("A10").Value = "Stop wage assignment"
$ws.Range("B10").Value = "https://www.illinoislegalaid.org/legal-information/stop-wage-assignment-letter"

# Add hyperlink for the new URL cell, matching existing pattern in the sheet
$ws.Hyperlinks.Add($ws.Range("B10"), "https://www.illinoislegalaid.org/legal-information/stop-wage-assignment-letter")

# Apply the same style as the other hyperlink cells in column B
$ws.Range("B10").Style = $ws.Range("B9").Style

# Update the selection to reflect final workbook state
$ws.Range("B14").Select()
